# Nudge the header group ("组合 17", shape id 18) on slide 1 down slightly.
# Target EMU offset: x=-1 (unchanged), y=117427 (was 100084).
# PowerPoint COM shape positions are expressed in points (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$grp = $s.Shapes.Item(1)

$targetTopEmu = 117427
$targetTopPts = $targetTopEmu / 12700

$grp.Top = $targetTopPts
